# Generate Report for Handoff
#
# Adds a newly-handed-off file (f6082fb2-0d97-4116-a91a-03b13b26fcd2.md) as a
# new row to the Overview / zh-cn / de-de status sheets, mirroring the shape
# of the existing 198b9554-... row, and grows each sheet's table to include
# the new row.

$wb = $excel.ActiveWorkbook

$repoBlobBase   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/155a168378ca4773075b1cee62e5eb18d2b2f720/e2e/"
$newFile        = "f6082fb2-0d97-4116-a91a-03b13b26fcd2.md"
$newFileDisplay = "e2e\f6082fb2-0d97-4116-a91a-03b13b26fcd2.md"
$newFileUrl     = $repoBlobBase + $newFile

# Helper-ish values: a leading apostrophe forces Excel to store the literal
# text instead of auto-coercing look-alike values ("True"/"False" -> bool,
# "" -> no cell at all). Resetting the style back to "Normal" afterwards
# drops the quote-prefix flag so the cell ends up a plain text cell.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Overview sheet: new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-TextValue $wsOverview.Range("A3") $newFile
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newFileUrl, "", "", $newFileDisplay)
$wsOverview.Range("B3").Style = "HyperLink"
Set-TextValue $wsOverview.Range("C3") ".md"
Set-TextValue $wsOverview.Range("D3") ""
Set-TextValue $wsOverview.Range("E3") "Ready for handoff"
Set-TextValue $wsOverview.Range("F3") "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-18 02:34:13"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# zh-cn sheet: new row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newFileUrl, "", "", $newFile)
$wsZhCn.Range("A3").Style = "HyperLink"
Set-TextValue $wsZhCn.Range("B3") ".md"
Set-TextValue $wsZhCn.Range("C3") "Ready for handoff"
Set-TextValue $wsZhCn.Range("D3") "e2e"
Set-TextValue $wsZhCn.Range("E3") "ht"
Set-TextValue $wsZhCn.Range("F3") "False"
Set-TextValue $wsZhCn.Range("G3") "f6082fb2-0d97-4116-a91a-03b13b26fcd2.dba96e427cedb1348ddbd22fbc04495a21b91893.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-18 02:34:08"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextValue $wsZhCn.Range("I3") ""
Set-TextValue $wsZhCn.Range("J3") ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextValue $wsZhCn.Range("L3") ""
Set-TextValue $wsZhCn.Range("M3") "True"
Set-TextValue $wsZhCn.Range("N3") ""
Set-TextValue $wsZhCn.Range("O3") "False"
Set-TextValue $wsZhCn.Range("P3") ""

$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------------
# de-de sheet: new row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newFileUrl, "", "", $newFile)
$wsDeDe.Range("A3").Style = "HyperLink"
Set-TextValue $wsDeDe.Range("B3") ".md"
Set-TextValue $wsDeDe.Range("C3") "Ready for handoff"
Set-TextValue $wsDeDe.Range("D3") "e2e"
Set-TextValue $wsDeDe.Range("E3") "ht"
Set-TextValue $wsDeDe.Range("F3") "False"
Set-TextValue $wsDeDe.Range("G3") "f6082fb2-0d97-4116-a91a-03b13b26fcd2.dba96e427cedb1348ddbd22fbc04495a21b91893.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-18 02:34:13"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextValue $wsDeDe.Range("I3") ""
Set-TextValue $wsDeDe.Range("J3") ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextValue $wsDeDe.Range("L3") ""
Set-TextValue $wsDeDe.Range("M3") "True"
Set-TextValue $wsDeDe.Range("N3") ""
Set-TextValue $wsDeDe.Range("O3") "False"
Set-TextValue $wsDeDe.Range("P3") ""

$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:P3"))

Write-Output "Generate Report for Handoff: added row for $newFile"
